# Update countries & provincias Spain
# Refresh the "Pais" dashboard sheet: update the "last updated" timestamp and
# refresh the per-country case counters. Several countries also changed row
# position because the source data is kept sorted by total cases, so both the
# country name (column A) and its statistics (columns B:H) are rewritten for
# every affected row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: "Datos actualizados a 1 de Abril de 2020 a las HH:MM" -------
$ws.Range("A1").Value = "Datos actualizados a 1 de Abril de 2020 a las 10:20"

# --- Per-row refresh: Country, Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes -------------------
$rows = @(
    @{ Row = 16;  Pais = "Austria";       B = 10366; C = 186; D = 1095; E = 9143; F = 198; G = 0; H = 128 },
    @{ Row = 22;  Pais = "Australia";     B = 4862;  C = 99;  D = 345;  E = 4496; F = 50;  G = 1; H = 21 },
    @{ Row = 28;  Pais = "Rusia";         B = 2777;  C = 440; D = 190;  E = 2563; F = 8;   G = 7; H = 24 },
    @{ Row = 29;  Pais = "Malasia";       B = 2766;  C = 0;   D = 537;  E = 2186; F = 94;  G = 0; H = 43 },
    @{ Row = 30;  Pais = "Chile";         B = 2738;  C = 0;   D = 156;  E = 2570; F = 14;  G = 0; H = 12 },
    @{ Row = 31;  Pais = "Polonia";       B = 2347;  C = 36;  D = 7;    E = 2305; F = 50;  G = 2; H = 35 },
    @{ Row = 32;  Pais = "Filipinas";     B = 2311;  C = 227; D = 50;   E = 2165; F = 1;   G = 8; H = 96 },
    @{ Row = 33;  Pais = "Ecuador";       B = 2302;  C = 0;   D = 58;   E = 2165; F = 100; G = 0; H = 79 },
    @{ Row = 34;  Pais = "Rumania";       B = 2245;  C = 0;   D = 220;  E = 1943; F = 62;  G = 0; H = 82 },
    @{ Row = 35;  Pais = "Luxemburgo";    B = 2178;  C = 0;   D = 80;   E = 2075; F = 31;  G = 0; H = 23 },
    @{ Row = 36;  Pais = "Japon";         B = 2178;  C = 0;   D = 424;  E = 1697; F = 69;  G = 0; H = 57 },
    @{ Row = 85;  Pais = "Kuwait";        B = 289;   C = 0;   D = 80;   E = 209;  F = 13;  G = 0; H = 0 },
    @{ Row = 122; Pais = "Paraguay";      B = 69;    C = 4;   D = 1;    E = 65;   F = 3;   G = 0; H = 3 },
    @{ Row = 123; Pais = "Gibraltar";     B = 69;    C = 0;   D = 34;   E = 35;   F = 0;   G = 0; H = 0 },
    @{ Row = 124; Pais = "Liechtenstein"; B = 68;    C = 0;   D = 0;    E = 68;   F = 0;   G = 0; H = 0 },
    @{ Row = 146; Pais = "Etiopia";       B = 29;    C = 3;   D = 4;    E = 25;   F = 2;   G = 0; H = 0 },
    @{ Row = 147; Pais = "Mali";          B = 28;    C = 0;   D = 0;    E = 26;   F = 0;   G = 0; H = 2 }
)

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.Pais
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 7).Value = $item.G
    $ws.Cells.Item($r, 8).Value = $item.H
}
